$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Simple price/volume updates in column D
Set-TextValue "D2" "247.66"
Set-TextValue "D7" "6.313"
Set-TextValue "D8" "0.8064"
Set-TextValue "D9" "0.8813"
Set-TextValue "D10" "0.1414"
Set-TextValue "D11" "0.07416"
Set-TextValue "D12" "0.03056"
Set-TextValue "D13" "0.03078"
Set-TextValue "D14" "0.09404"
Set-TextValue "D15" "3.883"
Set-TextValue "D16" "0.001570"
Set-TextValue "D17" "0.04790"
Set-TextValue "D18" "0.0005843"
Set-TextValue "D19" "0.006417"
Set-TextValue "D21" "0.0009966"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.690"
Set-TextValue "D24" "2.195"
Set-TextValue "D25" "0.3280"
Set-TextValue "D26" "0.1352"
Set-TextValue "D27" "0.01827"
Set-TextValue "D40" "0.03950"

# Rows 41-43 get rotated: KickToken -> row41, BKEXToken -> row42, CEJI -> row43
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006810"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003200"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.008446"
Set-TextValue "D45" "0.00005587"
Set-TextValue "D47" "0.4503"
Set-TextValue "D48" "0.2021"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
